# Weekly crime data update (cs-en-us-020pct.xlsx)
# - bumps the report "Number" and week-covering dates by one week
# - refreshes the Crime Complaints table (rows 15-29) with the new week's figures,
#   including a few cells that flip between a numeric count and the "0"/"***.*"
#   placeholder text used when there is no data / an undefined percent change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text: Volume/Number and the "week covering" date range.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 29   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/24/2022  Through  10/30/2022"

# ---------------------------------------------------------------------------
# 2. Plain numeric updates: the cell keeps its existing style/number format,
#    only the number itself changes.
# ---------------------------------------------------------------------------
$simpleUpdates = @{
    "D16" = "7"
    "E16" = "-71.428571428571"
    "F16" = "8"
    "G16" = "18"
    "H16" = "-55.555555555555"
    "I16" = "97"
    "J16" = "92"
    "K16" = "5.434782608695"
    "L16" = "27.631578947368"
    "M16" = "27.631578947368"
    "N16" = "-80.943025540275"
    "C17" = "3"
    "D17" = "1"
    "E17" = "200"
    "F17" = "9"
    "G17" = "14"
    "H17" = "-35.714285714285"
    "I17" = "88"
    "J17" = "60"
    "K17" = "46.666666666666"
    "L17" = "95.555555555555"
    "M17" = "79.591836734693"
    "N17" = "7.317073170731"
    "F18" = "2"
    "H18" = "-75"
    "L18" = "-4.838709677419"
    "M18" = "29.670329670329"
    "N18" = "-84.453227931488"
    "C19" = "21"
    "D19" = "32"
    "E19" = "-34.375"
    "F19" = "68"
    "G19" = "92"
    "H19" = "-26.086956521739"
    "I19" = "698"
    "J19" = "607"
    "K19" = "14.99176276771"
    "L19" = "71.921182266009"
    "M19" = "27.372262773722"
    "N19" = "-54.527687296416"
    "E20" = "-100"
    "F20" = "6"
    "G20" = "10"
    "H20" = "-40"
    "J20" = "79"
    "K20" = "-17.721518987341"
    "L20" = "91.176470588235"
    "N20" = "-93.305870236869"
    "D21" = "41"
    "E21" = "-36.585365853658"
    "F21" = "93"
    "G21" = "142"
    "H21" = "-34.507042253521"
    "I21" = "1078"
    "J21" = "912"
    "K21" = "18.201754385964"
    "L21" = "56.231884057971"
    "M21" = "34.413965087281"
    "N21" = "-72.18782249742"
    "F22" = "2"
    "H22" = "100"
    "F23" = "1"
    "G23" = "3"
    "H23" = "-66.666666666666"
    "I23" = "23"
    "K23" = "15"
    "L23" = "43.75"
    "M23" = "0"
    "C24" = "19"
    "D24" = "18"
    "E24" = "5.555555555555"
    "F24" = "68"
    "G24" = "91"
    "H24" = "-25.274725274725"
    "I24" = "1063"
    "J24" = "1076"
    "K24" = "-1.208178438661"
    "L24" = "0.472589792060"
    "M24" = "18.506131549609"
    "D25" = "4"
    "E25" = "25"
    "F25" = "18"
    "G25" = "16"
    "H25" = "12.5"
    "I25" = "178"
    "J25" = "151"
    "K25" = "17.880794701986"
    "L25" = "56.140350877193"
    "M25" = "-21.929824561403"
    "H26" = "-100"
    "G27" = "3"
    "H27" = "0"
    "G28" = "2"
    "J28" = "3"
    "K28" = "-66.666666666666"
    "G29" = "2"
    "J29" = "3"
    "K29" = "-66.666666666666"
}

foreach ($addr in $simpleUpdates.Keys) {
    $ws.Range($addr).Value = [double]$simpleUpdates[$addr]
}

# ---------------------------------------------------------------------------
# 3. Cells that switch from a number to the text placeholder "0" (means no
#    reported incidents) or "***.*" (undefined percent change), or vice versa.
#    A donor cell that already carries the desired style is copied (formats
#    only) onto the target after the value is written, so the target ends up
#    with exactly the same look as the other cells that use that style.
# ---------------------------------------------------------------------------

function Set-PlaceholderText($addr, $text, $donor) {
    # Force the cell to Text format so Excel stores the value as a string
    # instead of re-interpreting "0" as the number zero.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

function Set-NumberValue($addr, $value, $donor) {
    $ws.Range($addr).Value = [double]$value
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Donor cells: stable cells elsewhere in the table that are not touched by
# this edit and already have the style we want to replicate.
$textDonor = "G15"     # s=14 (right aligned, General) text cell
$intDonor  = "G18"     # s=16 (#,##0) integer cell
$pctDonor  = "K18"     # s=15 (#,##0.0) percent-style cell

# Number -> text ("0" / "***.*")
Set-PlaceholderText "F15" "0" $textDonor
Set-PlaceholderText "C18" "0" $textDonor
Set-PlaceholderText "D18" "0" $textDonor
Set-PlaceholderText "E18" "***.*" $textDonor
Set-PlaceholderText "C20" "0" $textDonor
Set-PlaceholderText "C22" "0" $textDonor
Set-PlaceholderText "D23" "0" $textDonor
Set-PlaceholderText "E23" "***.*" $textDonor
Set-PlaceholderText "F26" "0" $textDonor
Set-PlaceholderText "D27" "0" $textDonor
Set-PlaceholderText "E27" "***.*" $textDonor

# Text -> number
Set-NumberValue "C16" 2 $intDonor
Set-NumberValue "C23" 1 $intDonor
Set-NumberValue "D28" 1 $intDonor
Set-NumberValue "E28" -100 $pctDonor
Set-NumberValue "D29" 1 $intDonor
Set-NumberValue "E29" -100 $pctDonor

Write-Output "Done applying weekly crime data update."
